# Weekly refresh of the Fruta/Hortaliza data: the report window has rolled
# forward by one week, so the per-market rows are rotated to reflect the
# new weekly figures. Columns A-C, E-K (market/product identification)
# stay identical across rows, so we only need to move columns D, L-T.
#
# Mapping of new-row -> old-row that supplied its (D, L:T) values:
#   2 <- 8
#   3 <- 6
#   4 <- 2
#   5 <- 3
#   6 <- 5
#   7 <- 7 (unchanged)
#   8 <- 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the mutable columns (D and L:T) for each source row before any
# writes happen, so the rotation doesn't clobber values we still need.
# (Value2 is used instead of Value because it reliably returns/accepts the
# raw scalar instead of a bound property-accessor object.)
$snapshot = @{}
foreach ($r in 2..8) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

$rowMap = @{ 2 = 8; 3 = 6; 4 = 2; 5 = 3; 6 = 5; 7 = 7; 8 = 4 }

foreach ($destRow in 2..8) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2 = $data.D
    $ws.Cells.Item($destRow, 12).Value2 = $data.L
    $ws.Cells.Item($destRow, 13).Value2 = $data.M
    $ws.Cells.Item($destRow, 14).Value2 = $data.N
    $ws.Cells.Item($destRow, 15).Value2 = $data.O
    $ws.Cells.Item($destRow, 16).Value2 = $data.P
    $ws.Cells.Item($destRow, 17).Value2 = $data.Q
    $ws.Cells.Item($destRow, 18).Value2 = $data.R
    $ws.Cells.Item($destRow, 19).Value2 = $data.S
    $ws.Cells.Item($destRow, 20).Value2 = $data.T
}
